# Iteration 2 section edits
#
# 1. Move the "_GoBack" bookmark from the end of the "GUI will consist..."
#    paragraph to right after "CSC 450".
# 2. Split the "Once the hashmap is complete ... alphabetically." run into
#    several runs, wrapping "hashmap"/"hashmap"/"treemap" with spell-check
#    proofErr markers.
# 3. Split the "GUI will consist ... read in. " run, extending the sentence
#    and wrapping "gui" with spell-check proofErr markers.

$d = $word.ActiveDocument

function Insert-WordXmlFragment($range, [string]$fragment) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $fragment + '</w:p></w:body>' +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1. Relocate the "_GoBack" bookmark -------------------------------

$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$csc = $d.Content
$csc.Find.Execute("CSC 450", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cscEnd = $csc.End

# Temporarily extend the paragraph so the collapsed insertion point is no
# longer the paragraph's trailing mark, add the bookmark, then remove the
# placeholder again.
$tmp = $d.Range($cscEnd, $cscEnd)
$tmp.InsertAfter("#")
$bmSpot = $d.Range($cscEnd, $cscEnd)
$d.Bookmarks.Add("_GoBack", $bmSpot)
$d.Range($cscEnd, $cscEnd + 1).Delete()

# --- 2. Rewrite the hashmap/treemap sentence ---------------------------

$hm = $d.Content
$hm.Find.Execute(" Once the hashmap is complete I will feed the hashmap out to a treemap, which will order the output alphabetically.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hmRange = $d.Range($hm.Start, $hm.End)

$hmFrag = '<w:r><w:t xml:space="preserve"> Once the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>hashmap</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> is complete I will feed the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>hashmap</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> out to a </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>treemap</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, which will order the output alphabetically.</w:t></w:r>'

Insert-WordXmlFragment $hmRange $hmFrag

# --- 3. Rewrite the GUI sentence ---------------------------------------

$gui = $d.Content
$gui.Find.Execute("GUI will consist of a button which will allow user to choose the text file the would like to read in. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$guiRange = $d.Range($gui.Start, $gui.End)

$guiFrag = '<w:r><w:t xml:space="preserve">GUI will consist of a button which will allow user to choose the text </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">file the would like to read in and have the word counter be tested on. I am still unsure if I will have that output be shown on the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>gui</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> or have the output file be put somewhere else. </w:t></w:r>'

Insert-WordXmlFragment $guiRange $guiFrag

Write-Output "edits applied"
